$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @(
    "uni1C191C37",
    "uni1C1B1C37",
    "uni1C001C371C25",
    "uni1C031C371C25",
    "uni1C1D1C371C25",
    "uni1C101C37",
    "uni1C211C37",
    "uni1C001C371C24",
    "uni1C001C371C251C24",
    "uni1C031C371C24",
    "uni1C131C371C25",
    "uni1C1D1C371C24",
    "uni1C1D1C371C251C24",
    "uniE000"
)

$startRow = 433
$startSort = 454

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("C$row").Value = $newNames[$i]
    $ws.Range("D$row").Value = $startSort + $i
    $ws.Range("D$row").NumberFormat = "0.00"
}

$ws.Range("G37").Select()
